# Insert a new data row before the current row 34, shifting row 34..140
# down to 35..141 (dimension grows from A1:R140 to A1:R141), then
# populate the newly inserted row 34 with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlShiftDown = -4121
$ws.Rows.Item(34).Insert(-4121)

$ws.Cells.Item(34, 1).Value = 3
$ws.Cells.Item(34, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(34, 3).Value = "Coquimbo"
$ws.Cells.Item(34, 4).Value = 44607
$ws.Cells.Item(34, 5).Value = 5
$ws.Cells.Item(34, 6).Value = 100112052
$ws.Cells.Item(34, 7).Value = "Albahaca"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 148
$ws.Cells.Item(34, 11).Value = 3500
$ws.Cells.Item(34, 12).Value = 4000
$ws.Cells.Item(34, 13).Value = 3730
$ws.Cells.Item(34, 14).Value = "`$/docena de matas"
$ws.Cells.Item(34, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(34, 16).Value = 622
$ws.Cells.Item(34, 17).Value = 6
$ws.Cells.Item(34, 18).Value = "Hortaliza"
